$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.726.99'
$ws.Range("E2").Value = '  +0.89%  '

$ws.Range("D3").Value = '1.889.83'
$ws.Range("E3").Value = '  +0.77%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.51'
$ws.Range("E5").Value = '  +0.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4745'
$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2931'
$ws.Range("E8").Value = '  +0.74%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06532'
$ws.Range("E9").Value = '  +0.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.06'
$ws.Range("E10").Value = '  +0.58%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07814'
$ws.Range("E11").Value = '  +1.24%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.02'
$ws.Range("E12").Value = '  -0.62%  '

$ws.Range("D13").Value = '1.891.08'
$ws.Range("E13").Value = '  +0.79%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7358'
$ws.Range("E14").Value = '  -0.49%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.249'
$ws.Range("E15").Value = '  +2.42%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '284.17'
$ws.Range("E16").Value = '  +3.98%  '

$ws.Range("D17").Value = '30.728.36'
$ws.Range("E17").Value = '  +0.85%  '

$ws.Range("E18").Value = '  -1.57%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007536'
$ws.Range("E19").Value = '  -0.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  +0.04%  '

$ws.Range("D21").Value = '2.139.69'
$ws.Range("E21").Value = '  +0.51%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.333'
$ws.Range("E22").Value = '  +1.81%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.262'
$ws.Range("E24").Value = '  +1.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.236'
$ws.Range("E25").Value = '  -0.32%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.22'
$ws.Range("E26").Value = '  +0.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.92'
$ws.Range("E27").Value = '  +0.50%  '

$ws.Range("E28").Value = '  -0.41%  '

$ws.Range("E29").Value = '  -1.73%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09740'
$ws.Range("E30").Value = '  -3.28%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.499'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.303'
$ws.Range("E32").Value = '  -0.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.194'
$ws.Range("E33").Value = '  +2.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04862'
$ws.Range("E34").Value = '  +0.87%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.127'
$ws.Range("E35").Value = '  +0.24%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6983'
$ws.Range("E36").Value = '  -0.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.720'
$ws.Range("E37").Value = '  -0.12%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01904'
$ws.Range("E38").Value = '  +2.31%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.806'
$ws.Range("E39").Value = '  +1.98%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.395'
$ws.Range("E40").Value = '  +1.37%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '76.05'
$ws.Range("E41").Value = '  +7.21%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.019'
$ws.Range("E42").Value = '  +2.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4267'
$ws.Range("E43").Value = '  +1.60%  '

$ws.Range("E44").Value = '  -0.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8349'
$ws.Range("E45").Value = '  -0.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '101.55'
$ws.Range("E46").Value = '  -1.03%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.520'
$ws.Range("E47").Value = '  +2.07%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '35.68'
$ws.Range("E48").Value = '  +0.38%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.041'
$ws.Range("E49").Value = '  -0.06%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '918.94'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05753'
